# Update faturamento_diario.xlsx:
#  - Insert a new daily record (day 9, 07/2025) as a new row 8, pushing the
#    existing rows 8..69 down to 9..70.
#  - Correct three existing 07/2025 daily totals (rows 2, 6, 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 8 (shifts old rows 8..69 -> 9..70).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new daily record.
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = 15367.22
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 2025
$ws.Range("E8").Value = "07/2025"

# Correct previously-reported totals for early July 2025.
$ws.Range("B2").Value = 17999.03
$ws.Range("B6").Value = 15315.31
$ws.Range("B7").Value = 19905.55
